$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 6 (ano = 2025) metrics: total_customers and new_customers increased by 1,
# recompute dependent rate columns (new_rate, returning_rate)
$ws.Range("C6").Value = 434
$ws.Range("E6").Value = 125
$ws.Range("G6").Value = 28.80184331797235
$ws.Range("H6").Value = 71.19815668202764
